# "maj cas de charge" - Bump 3g et Right turn 2.2g
$wb = $excel.ActiveWorkbook

# --- Update data on "RIGHT TURN 2,2G" sheet (G/H/I columns, rows 6-47) ---
$ws = $wb.Worksheets.Item("RIGHT TURN 2,2G")

$ws.Range("G6").Value  = 1104.7070000000001
$ws.Range("H6").Value  = 1411.818
$ws.Range("I6").Value  = 581.02800000000002

$ws.Range("G7").Value  = -1104.06
$ws.Range("H7").Value  = 2116.1149999999998
$ws.Range("I7").Value  = 871.93100000000004

$ws.Range("G8").Value  = -149.577
$ws.Range("H8").Value  = 467.428
$ws.Range("I8").Value  = 3.7389999999999999

$ws.Range("G9").Value  = 148.98599999999999
$ws.Range("H9").Value  = 310.38799999999998
$ws.Range("I9").Value  = 2.4830000000000001

$ws.Range("G10").Value = 1104.7070000000001
$ws.Range("H10").Value = 1411.818
$ws.Range("I10").Value = 581.02800000000002

$ws.Range("G11").Value = -1104.06
$ws.Range("H11").Value = 2116.1149999999998
$ws.Range("I11").Value = 871.93100000000004

$ws.Range("G12").Value = -149.577
$ws.Range("H12").Value = 467.428
$ws.Range("I12").Value = 3.7389999999999999

$ws.Range("G13").Value = 148.98599999999999
$ws.Range("H13").Value = 310.38799999999998
$ws.Range("I13").Value = 2.4830000000000001

$ws.Range("G14").Value = -0.056
$ws.Range("H14").Value = -0.185
$ws.Range("I14").Value = -0.016

$ws.Range("G15").Value = -0.056
$ws.Range("H15").Value = -0.185
$ws.Range("I15").Value = -0.016

$ws.Range("G16").Value = 0
$ws.Range("H16").Value = -4305.5630000000001
$ws.Range("I16").Value = 1352.155

$ws.Range("G17").Value = 0
$ws.Range("H17").Value = -4305.5630000000001
$ws.Range("I17").Value = 1352.155

$ws.Range("H18").Value = -3332.9589999999998
$ws.Range("I18").Value = -1037.4079999999999

$ws.Range("H19").Value = -3332.9589999999998
$ws.Range("I19").Value = -1037.4079999999999

$ws.Range("G20").Value = 0
$ws.Range("H20").Value = -972.60400000000004
$ws.Range("I20").Value = 2389.5619999999999

$ws.Range("I23").Value = 2288.89

$ws.Range("I24").Value = 522.27200000000005

$ws.Range("G26").Value = -862.90099999999995
$ws.Range("H26").Value = -2044.0039999999999
$ws.Range("I26").Value = 113.249

$ws.Range("G27").Value = 774.31500000000005
$ws.Range("H27").Value = -1771.066
$ws.Range("I27").Value = 101.55800000000001

$ws.Range("G28").Value = -106.26900000000001
$ws.Range("H28").Value = 275.50200000000001
$ws.Range("I28").Value = -30.449000000000002

$ws.Range("G29").Value = 196.18100000000001
$ws.Range("H29").Value = 530.06399999999996
$ws.Range("I29").Value = -58.582999999999998

$ws.Range("G30").Value = 196.18100000000001
$ws.Range("H30").Value = 530.06399999999996
$ws.Range("I30").Value = -58.582999999999998

$ws.Range("G31").Value = -106.26900000000001
$ws.Range("H31").Value = 275.50200000000001
$ws.Range("I31").Value = -30.449000000000002

$ws.Range("G32").Value = 774.31500000000005
$ws.Range("H32").Value = -1771.066
$ws.Range("I32").Value = 101.55800000000001

$ws.Range("G33").Value = -862.90099999999995
$ws.Range("H33").Value = -2044.0039999999999
$ws.Range("I33").Value = 113.249

$ws.Range("H34").Value = 2880.7939999999999
$ws.Range("I34").Value = 2765.6990000000001

$ws.Range("H35").Value = 2880.7939999999999
$ws.Range("I35").Value = 2765.6990000000001

$ws.Range("H36").Value = -421.73500000000001
$ws.Range("I36").Value = 3358.085

$ws.Range("H37").Value = -421.73500000000001
$ws.Range("I37").Value = 3358.085

$ws.Range("G38").Value = -1.3260000000000001
$ws.Range("H38").Value = 128.71100000000001
$ws.Range("I38").Value = -12.798

$ws.Range("G39").Value = -1.3260000000000001
$ws.Range("H39").Value = 128.71100000000001
$ws.Range("I39").Value = -12.798

$ws.Range("I43").Value = 496.59

$ws.Range("I44").Value = 2374.9209999999998

$ws.Range("H47").Value = 3302.529
$ws.Range("I47").Value = -592.38599999999997

# --- View/selection changes ---

# "BRAKING 1,9G" sheet loses the tab-selected flag and scrolls back to top
$wsBraking = $wb.Worksheets.Item("BRAKING 1,9G")
$wsBraking.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2

# "RIGHT TURN 1 G" sheet selection moves
$wsRight1G = $wb.Worksheets.Item("RIGHT TURN 1 G")
$wsRight1G.Activate()
$wsRight1G.Range("G3:I44").Select()

# "RIGHT TURN 2,2G" becomes the active / tab-selected sheet, scrolled and
# reselected
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("N36").Select()
